# Burndown Chart - Avatar Supply Visibility: update the daily burndown
# numbers, clear out the two trailing (not-yet-reported) days, and shrink
# the chart's plotted range to match the new data extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Point" (remaining work) values for 9/30-10/11 ---
$ws.Range("B3").Value  = 29
$ws.Range("B4").Value  = 24
$ws.Range("B5").Value  = 19
$ws.Range("B6").Value  = 14
$ws.Range("B7").Value  = 14
$ws.Range("B8").Value  = 14
$ws.Range("B9").Value  = 9
$ws.Range("B10").Value = 4
$ws.Range("B11").Value = 0

# A new day's "Un-Completed" figure came in.
$ws.Range("C4").Value = 34

# The last two days (10/12 and 10/13) have not happened yet - clear them
# back out so only A12/A13 keep their date-cell formatting.
$ws.Range("A12:B13").ClearContents()

# Match the author's last active cell when they saved the file.
$ws.Range("C6").Select()

# --- Re-point the chart series at the new (smaller) data range ---
$chart = $ws.ChartObjects(1).Chart

$s1 = $chart.SeriesCollection(1)
$s1.Formula = "=SERIES('20190930-20191011'!`$B`$1,'20190930-20191011'!`$A`$2:`$A`$11,'20190930-20191011'!`$B`$2:`$B`$11,1)"

$s2 = $chart.SeriesCollection(2)
$s2.Formula = "=SERIES('20190930-20191011'!`$C`$1,'20190930-20191011'!`$A`$2:`$A`$11,'20190930-20191011'!`$C`$2:`$C`$13,2)"
